$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "searchCategory"
$ws.Range("C4").Value = "org.openqa.selenium.TimeoutException: Expected condition failed: waiting for visibility of Proxy element for: DefaultElementLocator 'By.xpath: //android.widget.TextView[@resource-id='com.zopsmart.stg.scarlet:id/tv_page_name']' (tried for 20 second(s) with 500 milliseconds interval)"
